# Adapt the column header formatting to the respective input file names:
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# then expose the sheet's used range as an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
$fv2504Headers = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

# Columns A..J (1..10) hold the "_old" headers -> rename to "_FV2410".
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}

# Column K (11) is the "diff" column and is left untouched.

# Columns L..U (12..21) hold the "_new" headers -> rename to "_FV2504".
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# Turn the used range into an Excel Table (ListObject): adds the autofilter and
# the table part, using the just-renamed header row as the column names.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U85"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
